$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "sd_temp"

$ws.Range("F2").Value = 1.32377526505859
$ws.Range("F3").Value = 2.76825999770389
$ws.Range("F4").Value = 2.37789475782806
$ws.Range("F5").Value = 1.5937428671309
$ws.Range("F6").Value = 2.41283060156728
$ws.Range("F7").Value = 2.09433005295105
$ws.Range("F8").Value = 1.43911221106763
$ws.Range("F9").Value = 2.45748461547577
$ws.Range("F10").Value = 2.2321299087631
$ws.Range("F11").Value = 0
$ws.Range("F12").Value = 0.94205954720293
$ws.Range("F13").Value = 2.58915702725425
$ws.Range("F14").Value = 2.00414353579197
